$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.724109292030334
$ws.Range("B1").Value = 2.542062997817993
$ws.Range("C1").Value = 3.276809930801392
$ws.Range("D1").Value = 1.232200980186462
$ws.Range("E1").Value = 0.8107597827911377
